# "Created A Fish A Day card." — adds two new card rows ("Chat" and
# "A Fish A Day!") to the card list on Sheet1, widens the new notes
# column (J) that holds them, and updates the sheet's selection/scroll
# state to match where the author was last looking.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 42: new "Chat" skill card.
$ws.Range("A42").Value = "Chat"
$ws.Range("B42").Value = 1
$ws.Range("C42").Value = "Common"
$ws.Range("D42").Value = "Skill"
$ws.Range("E42").Value = "Passivity"
$ws.Range("F42").Value = "If an enemy does not intend to attack, apply !M! Passivity. Increase this card's Passivity by !theRose:SecondMagic! for this combat."

# Row 43: new "A Fish A Day!" skill card.
$ws.Range("A43").Value = "A Fish A Day!"
$ws.Range("B43").Value = 1
$ws.Range("C43").Value = "Uncommon"
$ws.Range("D43").Value = "Skill"
$ws.Range("E43").Value = "Food"
$ws.Range("F43").Value = "Gain !M! Artifact. At the end of this turn, lose !M! Artifact. Exhaust."

# New column J width (holds card-design notes alongside the table).
$ws.Columns.Item(10).ColumnWidth = 27.8

# Restore the selection/scroll position the author left the sheet at.
$ws.Range("D36").Select()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
